$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Ref="D2"; Val="41.220.33"},
    @{Ref="E2"; Val="  -1.89%  "},
    @{Ref="D3"; Val="2.182.87"},
    @{Ref="E3"; Val="  -2.45%  "},
    @{Ref="E4"; Val="  -0.14%  "},
    @{Ref="D5"; Val="249.47"},
    @{Ref="E5"; Val="  -0.16%  "},
    @{Ref="D6"; Val="0.608"},
    @{Ref="E6"; Val="  -3.86%  "},
    @{Ref="D7"; Val="66.52"},
    @{Ref="E7"; Val="  -7.74%  "},
    @{Ref="E8"; Val="  -0.03%  "},
    @{Ref="D9"; Val="0.573"},
    @{Ref="E9"; Val="  -4.15%  "},
    @{Ref="D10"; Val="59.03"},
    @{Ref="E10"; Val="  +1.63%  "},
    @{Ref="D11"; Val="36.47"},
    @{Ref="E11"; Val="  -11.96%  "},
    @{Ref="D12"; Val="0.0929"},
    @{Ref="E12"; Val="  -5.24%  "},
    @{Ref="D13"; Val="0.104"},
    @{Ref="E13"; Val="  -1.22%  "},
    @{Ref="D14"; Val="6.91"},
    @{Ref="E14"; Val="  -3.98%  "},
    @{Ref="D15"; Val="2.505.36"},
    @{Ref="E15"; Val="  -2.48%  "},
    @{Ref="D16"; Val="14.38"},
    @{Ref="E16"; Val="  -4.77%  "},
    @{Ref="E17"; Val="  -1.78%  "},
    @{Ref="D18"; Val="2.133.04"},
    @{Ref="E18"; Val="  -4.31%  "},
    @{Ref="D19"; Val="41.162.47"},
    @{Ref="E19"; Val="  -1.97%  "},
    @{Ref="D20"; Val="0.0₃0946"},
    @{Ref="E20"; Val="  -3.26%  "},
    @{Ref="D21"; Val="71.68"},
    @{Ref="E21"; Val="  -2.22%  "},
    @{Ref="E22"; Val="  -2.82%  "},
    @{Ref="D23"; Val="230.65"},
    @{Ref="E23"; Val="  -2.51%  "},
    @{Ref="D24"; Val="2.04"},
    @{Ref="E24"; Val="  -5.66%  "},
    @{Ref="B25"; Val="WEMIXToken"},
    @{Ref="C25"; Val="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"},
    @{Ref="D25"; Val="3.77"},
    @{Ref="E25"; Val="  -5.79%  "},
    @{Ref="B26"; Val="Dai"},
    @{Ref="C26"; Val="https://coinranking.com/coin/MoTuySvg7+dai-dai"},
    @{Ref="D26"; Val="1.00"},
    @{Ref="E26"; Val="  +0.05%  "},
    @{Ref="D27"; Val="11.39"},
    @{Ref="E27"; Val="  +5.43%  "},
    @{Ref="D28"; Val="2.41"},
    @{Ref="E28"; Val="  -5.19%  "},
    @{Ref="E29"; Val="  -4.00%  "},
    @{Ref="B30"; Val="Monero"},
    @{Ref="C30"; Val="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"},
    @{Ref="D30"; Val="168.07"},
    @{Ref="E30"; Val="  -2.11%  "},
    @{Ref="B31"; Val="Toncoin"},
    @{Ref="C31"; Val="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"},
    @{Ref="D31"; Val="2.03"},
    @{Ref="E31"; Val="  -7.54%  "},
    @{Ref="D32"; Val="20.26"},
    @{Ref="E32"; Val="  -3.17%  "},
    @{Ref="E33"; Val="  -1.73%  "},
    @{Ref="D34"; Val="5.78"},
    @{Ref="E34"; Val="  +4.07%  "},
    @{Ref="D35"; Val="0.0736"},
    @{Ref="E35"; Val="  +0.31%  "},
    @{Ref="D36"; Val="0.122"},
    @{Ref="E36"; Val="  -3.25%  "},
    @{Ref="E37"; Val="  -4.65%  "},
    @{Ref="D38"; Val="3.99"},
    @{Ref="E38"; Val="  -0.93%  "},
    @{Ref="D39"; Val="24.64"},
    @{Ref="D40"; Val="0.0307"},
    @{Ref="E40"; Val="  +2.60%  "},
    @{Ref="D41"; Val="2.22"},
    @{Ref="E41"; Val="  -3.96%  "},
    @{Ref="D42"; Val="5.40"},
    @{Ref="E42"; Val="  +9.62%  "},
    @{Ref="D43"; Val="5.51"},
    @{Ref="E43"; Val="  -8.78%  "},
    @{Ref="B44"; Val="Celestia"},
    @{Ref="C44"; Val="https://coinranking.com/coin/YQcD0lBl7+celestia-tia"},
    @{Ref="D44"; Val="11.46"},
    @{Ref="E44"; Val="  -6.46%  "},
    @{Ref="B45"; Val="MultiversX"},
    @{Ref="C45"; Val="https://coinranking.com/coin/omwkOTglq+multiversx-egld"},
    @{Ref="D45"; Val="61.29"},
    @{Ref="E45"; Val="  -10.25%  "},
    @{Ref="D46"; Val="8.53"},
    @{Ref="E46"; Val="  -3.47%  "},
    @{Ref="D47"; Val="0.191"},
    @{Ref="E47"; Val="  -8.51%  "},
    @{Ref="D48"; Val="0.100"},
    @{Ref="E48"; Val="  -2.48%  "},
    @{Ref="E49"; Val="  -0.11%  "},
    @{Ref="E50"; Val="  -2.70%  "},
    @{Ref="E51"; Val="  -4.02%  "}
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Ref)
    $c.NumberFormat = "@"
    $c.Value = $u.Val
    $c.Style = "Normal"
}
